$wb = $excel.ActiveWorkbook

# --- "Sign up" sheet: update the username value (username113 -> username126) ---
$wsSignUp = $wb.Worksheets.Item("Sign up")
$wsSignUp.Range("I2").Value = "username126"

# --- "General Data" sheet: update the transaction date cell (D2) ---
# It moves from a numeric date (formatted w/ numFmtId 58) to a literal text
# value "4/28/2025" (numFmtId 49 / "@" text format), and the selection moves
# from E5 to D14.
$wsGeneral = $wb.Worksheets.Item("General Data")
$wsGeneral.Range("D2").NumberFormat = "@"
$wsGeneral.Range("D2").Value = "4/28/2025"
$wsGeneral.Range("D14").Select() | Out-Null

# --- Active tab moves from "General Data" to "Sign up" ---
# Select the General Data range first (above) so it does not end up being the
# active sheet; activating "Sign up" last makes it the workbook's active tab.
$wsSignUp.Select() | Out-Null
